# "funcionalidades de un libre y hojas"
#  - rename sheet "order" -> "Base de datos"
#  - add a new, second sheet named "Hoja1" after it
#  - give "Base de datos" an accent2/darker-25% tab color
#  - change the active selection on "Base de datos" to G9:H12 (active cell G9)

$wb = $excel.ActiveWorkbook

# --- first sheet: rename + tab color + selection -----------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Base de datos"

# Theme "Accent 2, Darker 25%" (theme index 5 / tint -0.249977111117893)
# resolved to RGB so it can be pushed through Tab.Color (an OLE BGR long).
$tabR = 0xC5
$tabG = 0x5A
$tabB = 0x11
$ws1.Tab.Color = $tabR + ($tabG * 256) + ($tabB * 65536)

[void]$ws1.Range("G9:H12").Select()

# --- second sheet: brand new blank sheet named "Hoja1" ------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Hoja1"

# Keep "Base de datos" as the selected/active tab (Add() activates the
# newly inserted sheet by default).
$ws1.Activate()
